# "added function to test all models"
# The test routine now also evaluates an SVM and a RandomForest model
# (in addition to the existing CNN/RNN/MLP models) against both the
# "audio files/" and "files/" datasets. This inserts two new result
# rows ahead of each existing 3-row block and refreshes the CNN/RNN/MLP
# numbers for the "files/" block (re-run results).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows in each block, working bottom-up so
# the row numbers used below stay valid.
$ws.Rows(5).Resize(2).Insert()
$ws.Rows(2).Resize(2).Insert()

# The inserted rows picked up the bold/bordered header formatting on
# B:E; reset those to the plain "Normal" style used by the rest of the
# data rows, then give column A the same bold/bordered/centered style
# used by the other A cells in the table.
$ws.Range("B2:E3").Style = "Normal"
$ws.Range("B7:E8").Style = "Normal"

$ws.Range("A4").Copy() | Out-Null
$ws.Range("A2:A3").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A7:A8").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# --- Block 1: "audio files/" ---------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "modelSVM.pkl"
$ws.Range("C2").Value = "audio files/"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = 0.9252336448598131

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "modelRandomForest.pkl"
$ws.Range("C3").Value = "audio files/"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = 0.9485981308411215

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "modelCNN.model"
$ws.Range("C4").Value = "audio files/"
$ws.Range("D4").Value = 0.4256819188594818
$ws.Range("E4").Value = 0.894859790802002

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "modelRNN.model"
$ws.Range("C5").Value = "audio files/"
$ws.Range("D5").Value = 0.1860886365175247
$ws.Range("E5").Value = 0.9439252614974976

$ws.Range("A6").Value = 0
$ws.Range("B6").Value = "modelMLP.model"
$ws.Range("C6").Value = "audio files/"
$ws.Range("D6").Value = 0.1019595563411713
$ws.Range("E6").Value = 0.9742990732192993

# --- Block 2: "files/" ----------------------------------------------
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "modelSVM.pkl"
$ws.Range("C7").Value = "files/"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 0.1212121212121212

$ws.Range("A8").Value = 0
$ws.Range("B8").Value = "modelRandomForest.pkl"
$ws.Range("C8").Value = "files/"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = 0.5454545454545454

$ws.Range("A9").Value = 0
$ws.Range("B9").Value = "modelCNN.model"
$ws.Range("C9").Value = "files/"
$ws.Range("D9").Value = 3.656139612197876
$ws.Range("E9").Value = 0.6363636255264282

$ws.Range("A10").Value = 0
$ws.Range("B10").Value = "modelRNN.model"
$ws.Range("C10").Value = "files/"
$ws.Range("D10").Value = 3.805042266845703
$ws.Range("E10").Value = 0.1818181872367859

$ws.Range("A11").Value = 0
$ws.Range("B11").Value = "modelMLP.model"
$ws.Range("C11").Value = "files/"
$ws.Range("D11").Value = 3.09204888343811
$ws.Range("E11").Value = 0.6060606241226196
